$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.131.43'
$ws.Range('E2').Value = '  -1.72%  '
$ws.Range('D3').Value = '2.407.46'
$ws.Range('E3').Value = '  -3.95%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.92%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.574'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.68%  '
$ws.Range('D9').Value = '2.406.22'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.105'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.11%  '
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.342'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.49'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.32%  '
$ws.Range('D15').Value = '2.841.84'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('E16').Value = '  -0.44%  '
$ws.Range('D17').Value = '61.481.17'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '2.410.53'
$ws.Range('E18').Value = '  -4.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.45%  '
$ws.Range('E20').Value = '  -4.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '315.04'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.37'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '0.0₃0952'
$ws.Range('E26').Value = '  -8.61%  '
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.529.34'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.69'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('E30').Value = '  -4.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '516.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.03'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.62%  '
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.85'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.88%  '
$ws.Range('E35').Value = '  -2.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.52'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.72%  '
$ws.Range('E38').Value = '  -5.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.375'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '139.57'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.34'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '140.81'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.55%  '
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.66'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0516'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.27%  '
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0924'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.50%  '
